$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values must remain plain text,
# matching the original inline-string cell type used for this price/volume table).
$updates = @{
    'D2' = '261.75'
    'E2' = '1.07%'
    'E3' = '1.11%'
    'D4' = '4.708'
    'E4' = '-0.17%'
    'D5' = '0.06205'
    'E5' = '2.79%'
    'D6' = '6.718'
    'E6' = '0.54%'
    'D7' = '0.8501'
    'E7' = '-0.96%'
    'D8' = '0.9171'
    'E8' = '-1.52%'
    'D9' = '0.1409'
    'E9' = '0.78%'
    'D10' = '0.04607'
    'E10' = '-7.70%'
    'D11' = '0.07089'
    'E11' = '0.41%'
    'D12' = '0.03144'
    'E12' = '0.39%'
    'D13' = '0.09059'
    'E13' = '-0.81%'
    'D14' = '0.001545'
    'E14' = '0.01%'
    'D15' = '0.0006153'
    'E15' = '1.51%'
    'D16' = '0.005991'
    'E16' = '-0.13%'
    'E17' = '0.22%'
    'E18' = '0.22%'
    'D19' = '2.180'
    'E19' = '0.64%'
    'E20' = '0.41%'
    'D21' = '0.1309'
    'E21' = '0.89%'
    'D22' = '4.089'
    'E22' = '-0.78%'
    'D23' = '0.04248'
    'E23' = '0.32%'
    'E24' = '-0.35%'
    'E25' = '-5.90%'
    'D26' = '0.0001200'
    'E26' = '0.09%'
    'E27' = '5.15%'
    'D40' = '0.03922'
    'E40' = '1.66%'
    'D41' = '0.1113'
    'E41' = '-0.46%'
    'D42' = '0.004133'
    'E42' = '5.05%'
    'E43' = '-4.78%'
    'D44' = '0.01392'
    'E44' = '-9.03%'
    'D45' = '0.00005172'
    'E45' = '0.93%'
    'E46' = '0.12%'
    'E47' = '-34.15%'
    'E48' = '27.94%'
    'E49' = '0.12%'
    'E50' = '0.12%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force a Text number format so the numeric-looking strings (prices, percentages)
    # are stored as literal text rather than being auto-converted to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
